$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper pattern used throughout this script:
#   1. Use Find/Replace (Replace:=2, wdReplaceAll) to swap the old text for
#      the new text. Word places the whole replacement into a single run
#      (taking on the formatting of the first run of the matched range).
#   2. Re-find the freshly inserted text to get its Range.
#   3. To split that merged run into two runs (when the diff calls for two
#      separate <w:r> elements) flip a character-formatting property (Bold)
#      on the sub-range that should become its own run: toggling it to $true
#      forces Word to carve out a new run (inheriting the full original
#      rPr, including rFonts/cs and szCs), then flipping it back to $false
#      removes the <w:b/> again while keeping the run boundary.
# ---------------------------------------------------------------------------

# --- Change 1: paragraph "c." ---------------------------------------------
# " Both failed and successful " + "projects" + " tend to follow the same trend"
# -> " " (kept) + "In general, there were more successful projects than failed ones." (new run)
$old1 = " Both failed and successful projects tend to follow the same trend"
$new1 = " In general, there were more successful projects than failed ones."
$r = $d.Content
$r.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

$r = $d.Content
$r.Find.Execute("In general, there were more successful projects than failed ones.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newRun = $d.Range($r.Start, $r.End)
$newRun.Font.Bold = $true
$newRun.Font.Bold = $false

# --- Change 2: paragraph "2." -----------------------------------------------
# "results." -> "results" (kept) + ", and the data set is not normally distributed." (new run)
$r = $d.Content
$r.Find.Execute("results.", $true, $false, $false, $false, $false, $true, 1, $false, "results, and the data set is not normally distributed.", 2)

$r = $d.Content
$r.Find.Execute("results, and the data set is not normally distributed.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitAt = $r.Start + 7   # length of "results"
$newRun = $d.Range($splitAt, $r.End)
$newRun.Font.Bold = $true
$newRun.Font.Bold = $false

# --- Change 3: paragraph "3." -----------------------------------------------
# "Some tables and graphs ... these " + "projects" + " campaigns for successful and unsuccessful ones."
# -> "Create a table/graph that excludes the outliers" (new run, "3. " run kept untouched)
$old3 = "Some tables and graphs that we could include are ones showing the most successful categories vs the most unsuccessful ones, and something showing the life of these projects campaigns for successful and unsuccessful ones."
$new3 = "Create a table/graph that excludes the outliers"
$r = $d.Content
$r.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

$r = $d.Content
$r.Find.Execute($new3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newRun = $d.Range($r.Start, $r.End)
$newRun.Font.Bold = $true
$newRun.Font.Bold = $false
